# Update "想去人数" (number of people interested) column F values
# on the "展览" (exhibition) and "全部类型" (all types) worksheets,
# reflecting the refreshed data snapshot generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 1052
$wsExhibit.Range("F4").Value  = 176
$wsExhibit.Range("F5").Value  = 2890
$wsExhibit.Range("F7").Value  = 271
$wsExhibit.Range("F8").Value  = 24
$wsExhibit.Range("F10").Value = 96
$wsExhibit.Range("F11").Value = 135
$wsExhibit.Range("F12").Value = 54
$wsExhibit.Range("F13").Value = 2717
$wsExhibit.Range("F14").Value = 974

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 1052
$wsAll.Range("F5").Value  = 176
$wsAll.Range("F6").Value  = 2890
$wsAll.Range("F8").Value  = 271
$wsAll.Range("F9").Value  = 24
$wsAll.Range("F12").Value = 96
$wsAll.Range("F13").Value = 135
$wsAll.Range("F14").Value = 54
$wsAll.Range("F15").Value = 2717
$wsAll.Range("F16").Value = 974
